$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 ---
$ws.Range("A2").Value = "Testing_ABC company"
$ws.Range("B2").Value = "Buyer"
# C2 (email / hyperlink) and D2 (TM_0001) stay the same.
# Address (F2) and Remarks (H2) get cleared out, GSTNumber moves from column F to column G.
$ws.Range("F2").ClearContents()
$ws.Range("H2").ClearContents()
$ws.Range("G2").Value = "hjbnm7845"
$ws.Range("J2").Value = $false

# --- Add new row 3 ---
$ws.Range("A3").Value = "Testing_ABC company1"
$ws.Range("B3").Value = "Both"
$ws.Range("C3").Value = "abc@gmail.com"
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:abc@gmail.com")
$ws.Range("C3").Style = "Hyperlink"
$ws.Range("D3").Value = "TM_0002"
$ws.Range("E3").Value = 1235485123
$ws.Range("G3").Value = 784512
$ws.Range("I3").Value = $true
$ws.Range("J3").Value = $false

# --- Extend the conditional-formatting ranges to cover the new row ---
$fcA = $ws.Range("A2").FormatConditions.Item(1)
$fcA.ModifyAppliesToRange($ws.Range("A2:A3"))
$fcD = $ws.Range("D2").FormatConditions.Item(1)
$fcD.ModifyAppliesToRange($ws.Range("D2:D4"))

# --- Selection ---
$null = $ws.Range("C15").Select()
